# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (G) with recalculated strike-count values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 5
    5  = 0
    6  = 5
    7  = 4
    8  = 3
    9  = 3
    10 = 7
    11 = 5
    12 = 4
    13 = 6
    14 = 3
    15 = 5
    16 = 2
    17 = 4
    18 = 2
    19 = 4
    20 = 3
    21 = 6
    22 = 3
    23 = 8
    24 = 10
    25 = 2
    26 = 3
    27 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
